# Auto-generated PowerShell-style Excel COM-interop script
# Applies updated 'F' column (想去人数 / want-to-go count) values across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 782
$ws.Cells.Item(3, 6).Value = 14439
$ws.Cells.Item(4, 6).Value = 14613
$ws.Cells.Item(5, 6).Value = 1373
$ws.Cells.Item(6, 6).Value = 1417
$ws.Cells.Item(7, 6).Value = 5950
$ws.Cells.Item(8, 6).Value = 998
$ws.Cells.Item(13, 6).Value = 1573
$ws.Cells.Item(14, 6).Value = 461
$ws.Cells.Item(16, 6).Value = 1237
$ws.Cells.Item(17, 6).Value = 1877
$ws.Cells.Item(19, 6).Value = 37
$ws.Cells.Item(20, 6).Value = 2305
$ws.Cells.Item(21, 6).Value = 579
$ws.Cells.Item(22, 6).Value = 837
$ws.Cells.Item(23, 6).Value = 3412
$ws.Cells.Item(25, 6).Value = 322
$ws.Cells.Item(26, 6).Value = 2479
$ws.Cells.Item(27, 6).Value = 618
$ws.Cells.Item(30, 6).Value = 1840
$ws.Cells.Item(31, 6).Value = 1090
$ws.Cells.Item(32, 6).Value = 1455
$ws.Cells.Item(33, 6).Value = 113
$ws.Cells.Item(34, 6).Value = 156
$ws.Cells.Item(35, 6).Value = 5028
$ws.Cells.Item(36, 6).Value = 4961
$ws.Cells.Item(38, 6).Value = 162
$ws.Cells.Item(39, 6).Value = 691
$ws.Cells.Item(40, 6).Value = 699
$ws.Cells.Item(41, 6).Value = 3325
$ws.Cells.Item(42, 6).Value = 48
$ws.Cells.Item(44, 6).Value = 347
$ws.Cells.Item(45, 6).Value = 126
$ws.Cells.Item(47, 6).Value = 4455
$ws.Cells.Item(48, 6).Value = 635
$ws.Cells.Item(49, 6).Value = 308

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 129
$ws.Cells.Item(20, 6).Value = 17

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 7736
$ws.Cells.Item(3, 6).Value = 267
$ws.Cells.Item(4, 6).Value = 927

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 7736
$ws.Cells.Item(3, 6).Value = 782
$ws.Cells.Item(4, 6).Value = 267
$ws.Cells.Item(5, 6).Value = 927
$ws.Cells.Item(7, 6).Value = 14439
$ws.Cells.Item(8, 6).Value = 14613
$ws.Cells.Item(9, 6).Value = 1373
$ws.Cells.Item(10, 6).Value = 1417
$ws.Cells.Item(11, 6).Value = 5950
$ws.Cells.Item(12, 6).Value = 998
$ws.Cells.Item(13, 6).Value = 129
$ws.Cells.Item(16, 6).Value = 1573
$ws.Cells.Item(17, 6).Value = 461
$ws.Cells.Item(18, 6).Value = 37
$ws.Cells.Item(19, 6).Value = 837
$ws.Cells.Item(20, 6).Value = 3412
$ws.Cells.Item(21, 6).Value = 322
$ws.Cells.Item(22, 6).Value = 2479
$ws.Cells.Item(23, 6).Value = 618
$ws.Cells.Item(25, 6).Value = 1840
$ws.Cells.Item(31, 6).Value = 1090
$ws.Cells.Item(32, 6).Value = 1455
$ws.Cells.Item(33, 6).Value = 156
$ws.Cells.Item(34, 6).Value = 17
$ws.Cells.Item(35, 6).Value = 5028
$ws.Cells.Item(36, 6).Value = 4961
$ws.Cells.Item(38, 6).Value = 691
$ws.Cells.Item(39, 6).Value = 3325
$ws.Cells.Item(41, 6).Value = 347
$ws.Cells.Item(42, 6).Value = 126
$ws.Cells.Item(45, 6).Value = 635
$ws.Cells.Item(46, 6).Value = 308
